$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「幸せは伝染性。喜びを広めよう」" (row 357) was removed entirely.
# Deleting the whole row shifts all subsequent rows up by one.
$ws.Rows.Item(357).Delete()
